# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, as published for the gh-pages data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 2465
$ws1.Range("F3").Value  = 758
$ws1.Range("F10").Value = 949
$ws1.Range("F12").Value = 132
$ws1.Range("F13").Value = 443
$ws1.Range("F17").Value = 24434
$ws1.Range("F18").Value = 2331
$ws1.Range("F19").Value = 151
$ws1.Range("F20").Value = 365
$ws1.Range("F22").Value = 75
$ws1.Range("F23").Value = 362
$ws1.Range("F25").Value = 81
$ws1.Range("F26").Value = 241
$ws1.Range("F28").Value = 75
$ws1.Range("F30").Value = 362
$ws1.Range("F32").Value = 447

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 273
$ws2.Range("F8").Value  = 140
$ws2.Range("F11").Value = 3646
$ws2.Range("F21").Value = 4130

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 171
$ws3.Range("F4").Value = 800

# Sheet "全部类型" (All types - combined listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 171
$ws4.Range("F4").Value  = 2465
$ws4.Range("F5").Value  = 800
$ws4.Range("F6").Value  = 758
$ws4.Range("F14").Value = 273
$ws4.Range("F18").Value = 949
$ws4.Range("F19").Value = 132
$ws4.Range("F20").Value = 443
$ws4.Range("F24").Value = 24434
$ws4.Range("F30").Value = 2331
$ws4.Range("F31").Value = 151
$ws4.Range("F33").Value = 365
$ws4.Range("F36").Value = 362
$ws4.Range("F38").Value = 241
$ws4.Range("F41").Value = 75
$ws4.Range("F46").Value = 447
$ws4.Range("F48").Value = 4130
